# Add two new student-course enrolment rows to the "student_courses" sheet
# and refresh the "No. of Students" summary count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("student_courses")

# Row 3: student 1 (Hasith fdfdf) enrolled in course 1 (Maths)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Hasith"
$ws.Range("C3").Value = "fdfdf"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Maths"

# Row 4: student 2 (Lee jkljk) enrolled in course 2 (Science)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Lee"
$ws.Range("C4").Value = "jkljk"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "Science"

# Keep the newly written cells unstyled, matching the rest of the data rows.
$ws.Range("A3:E4").Style = "Normal"

# Update the "No. of Students" summary cell (L4) to reflect the 3 enrolled rows.
$ws.Range("L4").Value = 3
